# st.experimental_rerun() deprecated hence changed to st.rerun()
# Appending newly-logged dedupe run rows (9-18) to the dashboard log sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 1).Value = 'Unique'
$ws.Cells.Item(9, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(9, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(9, 4).Value = 0.7465200424194336
$ws.Cells.Item(9, 5).Value = '16/07/2024 16:23:55'
$ws.Cells.Item(9, 6).Value = '16/07/2024 16:24:35'
$ws.Cells.Item(9, 7).Value = 40
$ws.Cells.Item(9, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(9, 9).Value = 2917
$ws.Cells.Item(9, 10).Value = 4
$ws.Cells.Item(9, 11).Value = 16

# Row 10
$ws.Cells.Item(10, 1).Value = 'All'
$ws.Cells.Item(10, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(10, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(10, 4).Value = 0.7465200424194336
$ws.Cells.Item(10, 5).Value = '16/07/2024 16:23:55'
$ws.Cells.Item(10, 6).Value = '16/07/2024 16:24:37'
$ws.Cells.Item(10, 7).Value = 43
$ws.Cells.Item(10, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(10, 9).Value = 2921
$ws.Cells.Item(10, 10).Value = 4
$ws.Cells.Item(10, 11).Value = 16

# Row 11
$ws.Cells.Item(11, 1).Value = 'Unique'
$ws.Cells.Item(11, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(11, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(11, 4).Value = 0.7465200424194336
$ws.Cells.Item(11, 5).Value = '16/07/2024 16:23:55'
$ws.Cells.Item(11, 6).Value = '16/07/2024 16:24:47'
$ws.Cells.Item(11, 7).Value = 52
$ws.Cells.Item(11, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(11, 9).Value = 2917
$ws.Cells.Item(11, 10).Value = 4
$ws.Cells.Item(11, 11).Value = 16

# Row 12
$ws.Cells.Item(12, 1).Value = 'Unique'
$ws.Cells.Item(12, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(12, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(12, 4).Value = 0.7465200424194336
$ws.Cells.Item(12, 5).Value = '16/07/2024 16:23:55'
$ws.Cells.Item(12, 6).Value = '16/07/2024 16:24:51'
$ws.Cells.Item(12, 7).Value = 56
$ws.Cells.Item(12, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(12, 9).Value = 2917
$ws.Cells.Item(12, 10).Value = 4
$ws.Cells.Item(12, 11).Value = 16

# Row 13
$ws.Cells.Item(13, 1).Value = 'All'
$ws.Cells.Item(13, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(13, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(13, 4).Value = 0.7465200424194336
$ws.Cells.Item(13, 5).Value = '16/07/2024 16:23:55'
$ws.Cells.Item(13, 6).Value = '16/07/2024 16:24:53'
$ws.Cells.Item(13, 7).Value = 59
$ws.Cells.Item(13, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(13, 9).Value = 2921
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(13, 11).Value = 16

# Row 14
$ws.Cells.Item(14, 1).Value = 'Unique'
$ws.Cells.Item(14, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(14, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(14, 4).Value = 0.7465200424194336
$ws.Cells.Item(14, 5).Value = '16/07/2024 18:11:03'
$ws.Cells.Item(14, 6).Value = '16/07/2024 18:11:45'
$ws.Cells.Item(14, 7).Value = 43
$ws.Cells.Item(14, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(14, 9).Value = 2917
$ws.Cells.Item(14, 10).Value = 4
$ws.Cells.Item(14, 11).Value = 16

# Row 15
$ws.Cells.Item(15, 1).Value = 'All'
$ws.Cells.Item(15, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(15, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(15, 4).Value = 0.7465200424194336
$ws.Cells.Item(15, 5).Value = '16/07/2024 18:11:03'
$ws.Cells.Item(15, 6).Value = '16/07/2024 18:11:48'
$ws.Cells.Item(15, 7).Value = 45
$ws.Cells.Item(15, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(15, 9).Value = 2921
$ws.Cells.Item(15, 10).Value = 4
$ws.Cells.Item(15, 11).Value = 16

# Row 16
$ws.Cells.Item(16, 1).Value = 'Unique'
$ws.Cells.Item(16, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(16, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(16, 4).Value = 0.7465200424194336
$ws.Cells.Item(16, 5).Value = '16/07/2024 18:11:03'
$ws.Cells.Item(16, 6).Value = '16/07/2024 18:12:58'
$ws.Cells.Item(16, 7).Value = 115
$ws.Cells.Item(16, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(16, 9).Value = 2917
$ws.Cells.Item(16, 10).Value = 4
$ws.Cells.Item(16, 11).Value = 16

# Row 17
$ws.Cells.Item(17, 1).Value = 'Unique'
$ws.Cells.Item(17, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(17, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(17, 4).Value = 0.7465200424194336
$ws.Cells.Item(17, 5).Value = '16/07/2024 18:11:03'
$ws.Cells.Item(17, 6).Value = '16/07/2024 18:13:01'
$ws.Cells.Item(17, 7).Value = 118
$ws.Cells.Item(17, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(17, 9).Value = 2917
$ws.Cells.Item(17, 10).Value = 4
$ws.Cells.Item(17, 11).Value = 16

# Row 18
$ws.Cells.Item(18, 1).Value = 'All'
$ws.Cells.Item(18, 2).Value = 'Akash Shahapure Test'
$ws.Cells.Item(18, 3).Value = 'akash.shahapure@haqdarshak.com'
$ws.Cells.Item(18, 4).Value = 0.7465200424194336
$ws.Cells.Item(18, 5).Value = '16/07/2024 18:11:03'
$ws.Cells.Item(18, 6).Value = '16/07/2024 18:13:05'
$ws.Cells.Item(18, 7).Value = 122
$ws.Cells.Item(18, 8).Value = 'cases_report_Shapoorji Pallonji Group_MH_(All States)_2024-07-16'
$ws.Cells.Item(18, 9).Value = 2921
$ws.Cells.Item(18, 10).Value = 4
$ws.Cells.Item(18, 11).Value = 16

Write-Host "Added rows 9-18 to sheet"
